$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(7, 8).Value = 56142.855
$ws.Cells.Item(7, 9).Value = 20000
$ws.Cells.Item(7, 10).Value = 59947.367
$ws.Cells.Item(7, 11).Value = 20000
$ws.Cells.Item(7, 12).Value = 59947.367
$ws.Cells.Item(7, 13).Value = -19888
$ws.Cells.Item(7, 14).Value = -60171.367
$ws.Cells.Item(14, 8).Value = 56142.855
$ws.Cells.Item(14, 9).Value = 20000
$ws.Cells.Item(14, 10).Value = 59947.367
$ws.Cells.Item(14, 11).Value = 20000
$ws.Cells.Item(14, 12).Value = 59947.367
$ws.Cells.Item(14, 13).Value = -19809
$ws.Cells.Item(14, 14).Value = -60329.367
$ws.Cells.Item(43, 8).Value = 3657
$ws.Cells.Item(43, 9).Value = 1638
$ws.Cells.Item(43, 10).Value = 5676
$ws.Cells.Item(43, 11).Value = 1638
$ws.Cells.Item(43, 12).Value = 5676
$ws.Cells.Item(43, 13).Value = -1569
$ws.Cells.Item(43, 14).Value = -5814
$ws.Cells.Item(70, 8).Value = 903.6
$ws.Cells.Item(70, 9).Value = 614.4
$ws.Cells.Item(70, 10).Value = 1000
$ws.Cells.Item(70, 11).Value = 1843.2
$ws.Cells.Item(70, 12).Value = 3000
$ws.Cells.Item(70, 13).Value = -1573.2
$ws.Cells.Item(70, 14).Value = -3540
$ws.Cells.Item(73, 8).Value = 903.6
$ws.Cells.Item(73, 9).Value = 614.4
$ws.Cells.Item(73, 10).Value = 1000
$ws.Cells.Item(73, 11).Value = 1843.2
$ws.Cells.Item(73, 12).Value = 3000
$ws.Cells.Item(73, 13).Value = -907.1999999999998
$ws.Cells.Item(73, 14).Value = -4872
$ws.Cells.Item(127, 8).Value = 1167.7059
$ws.Cells.Item(127, 9).Value = 896.2308
$ws.Cells.Item(127, 10).Value = 2050
$ws.Cells.Item(127, 11).Value = 2688.6924
$ws.Cells.Item(127, 12).Value = 6150
$ws.Cells.Item(127, 13).Value = 2271.3076
$ws.Cells.Item(127, 14).Value = -16070
$ws.Cells.Item(132, 8).Value = 1359.8422
$ws.Cells.Item(132, 9).Value = 915.74194
$ws.Cells.Item(132, 10).Value = 3326.5715
$ws.Cells.Item(132, 11).Value = 2747.22582
$ws.Cells.Item(132, 12).Value = 9979.7145
$ws.Cells.Item(132, 13).Value = -217.2258200000001
$ws.Cells.Item(132, 14).Value = -15039.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5803.2856
$ws.Cells.Item(32, 9).Value = 2579.9424
$ws.Cells.Item(32, 10).Value = 21040.908
$ws.Cells.Item(32, 11).Value = 2579.9424
$ws.Cells.Item(32, 12).Value = 21040.908
$ws.Cells.Item(32, 13).Value = -2292.9424
$ws.Cells.Item(32, 14).Value = -21614.908
$ws.Cells.Item(110, 8).Value = 1241.1143
$ws.Cells.Item(110, 9).Value = 927.0645
$ws.Cells.Item(110, 10).Value = 3675
$ws.Cells.Item(110, 11).Value = 927.0645
$ws.Cells.Item(110, 12).Value = 3675
$ws.Cells.Item(110, 13).Value = 1117.9355
$ws.Cells.Item(110, 14).Value = -7765
$ws.Cells.Item(122, 8).Value = 2561.4666
$ws.Cells.Item(122, 9).Value = 2583.2273
$ws.Cells.Item(122, 10).Value = 2501.625
$ws.Cells.Item(122, 11).Value = 7749.6819
$ws.Cells.Item(122, 12).Value = 7504.875
$ws.Cells.Item(122, 13).Value = -5299.6819
$ws.Cells.Item(122, 14).Value = -12404.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 1138146.1
$ws.Cells.Item(105, 9).Value = 2274143.2
$ws.Cells.Item(105, 10).Value = 2149.1
$ws.Cells.Item(105, 11).Value = 2274143.2
$ws.Cells.Item(105, 12).Value = 2149.1
$ws.Cells.Item(105, 13).Value = -2272396.2
$ws.Cells.Item(105, 14).Value = -5643.1
$ws.Cells.Item(134, 8).Value = 1989.96
$ws.Cells.Item(134, 9).Value = 1845.1904
$ws.Cells.Item(134, 10).Value = 2750
$ws.Cells.Item(134, 11).Value = 5535.5712
$ws.Cells.Item(134, 12).Value = 8250
$ws.Cells.Item(134, 13).Value = -3000.5712
$ws.Cells.Item(134, 14).Value = -13320

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1294.1448
$ws.Cells.Item(31, 9).Value = 934.8182
$ws.Cells.Item(31, 10).Value = 1440.537
$ws.Cells.Item(31, 11).Value = 934.8182
$ws.Cells.Item(31, 12).Value = 1440.537
$ws.Cells.Item(31, 13).Value = -639.8182
$ws.Cells.Item(31, 14).Value = -2030.537
$ws.Cells.Item(33, 8).Value = 24565.75
$ws.Cells.Item(33, 9).Value = 4131.5
$ws.Cells.Item(33, 10).Value = 45000
$ws.Cells.Item(33, 11).Value = 4131.5
$ws.Cells.Item(33, 12).Value = 45000
$ws.Cells.Item(33, 13).Value = -3752.5
$ws.Cells.Item(33, 14).Value = -45758
$ws.Cells.Item(34, 8).Value = 1294.1448
$ws.Cells.Item(34, 9).Value = 934.8182
$ws.Cells.Item(34, 10).Value = 1440.537
$ws.Cells.Item(34, 11).Value = 934.8182
$ws.Cells.Item(34, 12).Value = 1440.537
$ws.Cells.Item(34, 13).Value = -732.8182
$ws.Cells.Item(34, 14).Value = -1844.537
$ws.Cells.Item(132, 8).Value = 2419.2222
$ws.Cells.Item(132, 9).Value = 2022.6086
$ws.Cells.Item(132, 10).Value = 4699.75
$ws.Cells.Item(132, 11).Value = 6067.825800000001
$ws.Cells.Item(132, 12).Value = 14099.25
$ws.Cells.Item(132, 13).Value = -3537.825800000001
$ws.Cells.Item(132, 14).Value = -19159.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(42, 8).Value = 4391.1665
$ws.Cells.Item(42, 9).Value = 201.5
$ws.Cells.Item(42, 10).Value = 6486
$ws.Cells.Item(42, 11).Value = 604.5
$ws.Cells.Item(42, 12).Value = 19458
$ws.Cells.Item(42, 13).Value = -70.5
$ws.Cells.Item(42, 14).Value = -20526
$ws.Cells.Item(47, 8).Value = 980
$ws.Cells.Item(47, 9).Value = 466.66666
$ws.Cells.Item(47, 10).Value = 1750
$ws.Cells.Item(47, 11).Value = 1399.99998
$ws.Cells.Item(47, 12).Value = 5250
$ws.Cells.Item(47, 13).Value = -968.9999800000001
$ws.Cells.Item(47, 14).Value = -6112
$ws.Cells.Item(107, 8).Value = 349.94595
$ws.Cells.Item(107, 9).Value = 218.64062
$ws.Cells.Item(107, 10).Value = 1190.3
$ws.Cells.Item(107, 11).Value = 655.92186
$ws.Cells.Item(107, 12).Value = 3570.9
$ws.Cells.Item(107, 13).Value = 1264.07814
$ws.Cells.Item(107, 14).Value = -7410.9
$ws.Cells.Item(137, 8).Value = 647865.25
$ws.Cells.Item(137, 9).Value = 945.2632
$ws.Cells.Item(137, 10).Value = 1672155.1
$ws.Cells.Item(137, 11).Value = 2835.7896
$ws.Cells.Item(137, 12).Value = 5016465.300000001
$ws.Cells.Item(137, 13).Value = 2264.2104
$ws.Cells.Item(137, 14).Value = -5026665.300000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2469.3684
$ws.Cells.Item(102, 9).Value = 1973.963
$ws.Cells.Item(102, 10).Value = 3685.3635
$ws.Cells.Item(102, 11).Value = 1973.963
$ws.Cells.Item(102, 12).Value = 3685.3635
$ws.Cells.Item(102, 13).Value = -351.963
$ws.Cells.Item(102, 14).Value = -6929.363499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 1750.25
$ws.Cells.Item(136, 9).Value = 926.8
$ws.Cells.Item(136, 10).Value = 3122.6667
$ws.Cells.Item(136, 11).Value = 2780.4
$ws.Cells.Item(136, 12).Value = 9368.000100000001
$ws.Cells.Item(136, 13).Value = -230.3999999999996
$ws.Cells.Item(136, 14).Value = -14468.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(44, 8).Value = 24500
$ws.Cells.Item(44, 10).Value = 24500
$ws.Cells.Item(44, 12).Value = 24500
$ws.Cells.Item(44, 14).Value = -25608
$ws.Cells.Item(122, 8).Value = 1945.6364
$ws.Cells.Item(122, 9).Value = 1641.3077
$ws.Cells.Item(122, 10).Value = 3076
$ws.Cells.Item(122, 11).Value = 4923.9231
$ws.Cells.Item(122, 12).Value = 9228
$ws.Cells.Item(122, 13).Value = -2473.9231
$ws.Cells.Item(122, 14).Value = -14128
